$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.948.28'
$ws.Range("E2").Value = '  -1.79%  '
$ws.Range("D3").Value = '2.378.76'
$ws.Range("E3").Value = '  +3.36%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '300.55'
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").Value = '98.48'
$ws.Range("E6").Value = '  -2.79%  '
$ws.Range("D7").Value = '0.566'
$ws.Range("E7").Value = '  -0.46%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").Value = '0.508'
$ws.Range("E9").Value = '  -4.21%  '
$ws.Range("D10").Value = '34.46'
$ws.Range("E10").Value = '  -6.16%  '
$ws.Range("D11").Value = '0.0787'
$ws.Range("E11").Value = '  -2.06%  '
$ws.Range("D12").Value = '7.11'
$ws.Range("E12").Value = '  -4.85%  '
$ws.Range("E13").Value = '  -0.37%  '
$ws.Range("D14").Value = '2.743.71'
$ws.Range("E14").Value = '  +3.48%  '
$ws.Range("D15").Value = '2.388.07'
$ws.Range("E15").Value = '  +3.66%  '
$ws.Range("D16").Value = '0.824'
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("D17").Value = '13.75'
$ws.Range("E17").Value = '  -1.77%  '
$ws.Range("D18").Value = '45.925.91'
$ws.Range("E18").Value = '  -1.79%  '
$ws.Range("D19").Value = '12.64'
$ws.Range("E19").Value = '  -6.67%  '
$ws.Range("D20").Value = '0.0₃0946'
$ws.Range("E20").Value = '  -0.66%  '
$ws.Range("D21").Value = '6.04'
$ws.Range("E21").Value = '  -1.28%  '
$ws.Range("D22").Value = '66.69'
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").Value = '243.18'
$ws.Range("E23").Value = '  -2.23%  '
$ws.Range("E24").Value = '  -5.35%  '
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("D26").Value = '1.92'
$ws.Range("E26").Value = '  -2.27%  '
$ws.Range("D27").Value = '39.48'
$ws.Range("E27").Value = '  -10.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.20'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.31%  '
$ws.Range("D29").Value = '9.73'
$ws.Range("E29").Value = '  -2.24%  '
$ws.Range("D30").Value = '20.93'
$ws.Range("E30").Value = '  +3.55%  '
$ws.Range("D31").Value = '3.73'
$ws.Range("E31").Value = '  +17.40%  '
$ws.Range("E32").Value = '  +6.93%  '
$ws.Range("D33").Value = '5.52'
$ws.Range("E33").Value = '  -4.73%  '
$ws.Range("D34").Value = '146.57'
$ws.Range("E34").Value = '  -0.46%  '
$ws.Range("D35").Value = '0.0771'
$ws.Range("E35").Value = '  -3.57%  '
$ws.Range("E36").Value = '  -1.49%  '
$ws.Range("E37").Value = '  +6.27%  '
$ws.Range("E38").Value = '  -2.63%  '
$ws.Range("D39").Value = '14.94'
$ws.Range("E39").Value = '  -7.39%  '
$ws.Range("D40").Value = '3.87'
$ws.Range("E40").Value = '  -3.69%  '
$ws.Range("E41").Value = '  -2.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.20'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.29%  '
$ws.Range("D43").Value = '1.941.01'
$ws.Range("E43").Value = '  +3.99%  '
$ws.Range("D45").Value = '92.13'
$ws.Range("E45").Value = '  +4.64%  '
$ws.Range("E46").Value = '  -9.33%  '
$ws.Range("D47").Value = '8.49'
$ws.Range("E47").Value = '  +5.12%  '
$ws.Range("D48").Value = '0.185'
$ws.Range("E48").Value = '  -6.21%  '
$ws.Range("D49").Value = '98.96'
$ws.Range("E49").Value = '  +1.96%  '
$ws.Range("D50").Value = '2.613.90'
$ws.Range("E50").Value = '  +3.41%  '
$ws.Range("D51").Value = '68.74'
$ws.Range("E51").Value = '  -7.71%  '
